# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Each row's Price (column D) and Volume(1h) (column E) cell is an inline
# text string in the source workbook (e.g. "67.630.20", "  +1.13%  "),
# never a real number, so every write below must land as text. Plain
# numeric-looking strings ("599.21", "7.00", "1.00", ...) would
# otherwise be auto-coerced to numbers by the COM value setter, silently
# dropping trailing zeros / exact formatting, so those cells are forced to
# text via a temporary "@" (Text) number format and then restored to the
# workbook's default "Normal" style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.685.04'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.498.06'
$ws.Range('E3').Value = '  -0.06%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '599.21'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '180.42'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  +4.43%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.618'
$r.Style = 'Normal'
$ws.Range('E7').Value = '  +5.91%  '
$ws.Range('D9').Value = '3.498.49'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +5.46%  '
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '7.00'
$r.Style = 'Normal'
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '4.107.38'
$ws.Range('E13').Value = '  +0.10%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '32.31'
$r.Style = 'Normal'
$ws.Range('E14').Value = '  +10.45%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '0.135'
$r.Style = 'Normal'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '67.642.23'
$ws.Range('E16').Value = '  +1.06%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '0.0000179'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '3.498.90'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('E19').Value = '  +1.06%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '14.32'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  +0.27%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '392.41'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '7.97'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +0.33%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '73.22'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('E25').Value = '  +0.15%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '5.76'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('E27').Value = '  +1.64%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '10.39'
$r.Style = 'Normal'
$ws.Range('E28').Value = '  +2.14%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '0.176'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -2.20%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('E34').Value = '  -0.40%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '7.45'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  +1.04%  '
$ws.Range('E37').Value = '  +0.52%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '162.98'
$r.Style = 'Normal'
$ws.Range('E38').Value = '  -0.87%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.887'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  +1.10%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '2.84'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  +11.96%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '1.90'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -0.40%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '6.86'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('E43').Value = '  +0.38%  '
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '26.49'
$r.Style = 'Normal'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').Value = '2.839.63'
$ws.Range('E45').Value = '  +0.01%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '26.92'
$r.Style = 'Normal'
$ws.Range('E46').Value = '  -1.00%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '0.0727'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  -1.00%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '41.69'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('E49').Value = '  -0.54%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '337.82'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('E51').Value = '  -0.67%  '
